$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct value assignments for non-percentage text cells (dates, pressure, radiation, temperature)
# These keep their original style (s=3) and inline/shared string text type automatically.
$ws.Range("E2").Value = "2026-02-21 17:18:20"
$ws.Range("K2").Value = "12.7 MJ/m2"
$ws.Range("E3").Value = "2026-02-21 17:18:23"
$ws.Range("K3").Value = "16.1 MJ/m2"
$ws.Range("E4").Value = "2026-02-21 17:18:25"
$ws.Range("J4").Value = "1029.4 hPa"
$ws.Range("K4").Value = "14.6 MJ/m2"
$ws.Range("O4").Value = "9.8 °C"
$ws.Range("E5").Value = "2026-02-21 17:18:28"
$ws.Range("K5").Value = "16.0 MJ/m2"
$ws.Range("O5").Value = "3.4 °C"
$ws.Range("E6").Value = "2026-02-21 17:18:30"
$ws.Range("K6").Value = "15.2 MJ/m2"
$ws.Range("O6").Value = "10.3 °C"
$ws.Range("E7").Value = "2026-02-21 17:18:33"
$ws.Range("K7").Value = "15.4 MJ/m2"
$ws.Range("O7").Value = "13.4 °C"
$ws.Range("E8").Value = "2026-02-21 17:18:35"
$ws.Range("K8").Value = "16.0 MJ/m2"
$ws.Range("O8").Value = "10.7 °C"
$ws.Range("E9").Value = "2026-02-21 17:18:37"
$ws.Range("K9").Value = "15.1 MJ/m2"
$ws.Range("E10").Value = "2026-02-21 17:18:40"
$ws.Range("K10").Value = "14.9 MJ/m2"
$ws.Range("O10").Value = "8.7 °C"
$ws.Range("E11").Value = "2026-02-21 17:18:42"
$ws.Range("E12").Value = "2026-02-21 17:18:44"
$ws.Range("E13").Value = "2026-02-21 17:18:46"
$ws.Range("J13").Value = "1032.0 hPa"
$ws.Range("K13").Value = "16.0 MJ/m2"
$ws.Range("O13").Value = "4.9 °C"
$ws.Range("E14").Value = "2026-02-21 17:18:49"
$ws.Range("K14").Value = "15.3 MJ/m2"
$ws.Range("E15").Value = "2026-02-21 17:18:51"
$ws.Range("O15").Value = "14.1 °C"
$ws.Range("E16").Value = "2026-02-21 17:18:54"
$ws.Range("E17").Value = "2026-02-21 17:18:56"
$ws.Range("K17").Value = "16.9 MJ/m2"
$ws.Range("O17").Value = "8.5 °C"
$ws.Range("E18").Value = "2026-02-21 17:18:59"
$ws.Range("K18").Value = "15.3 MJ/m2"
$ws.Range("O18").Value = "8.5 °C"
$ws.Range("E19").Value = "2026-02-21 17:19:01"
$ws.Range("K19").Value = "15.1 MJ/m2"
$ws.Range("O19").Value = "7.8 °C"
$ws.Range("E20").Value = "2026-02-21 17:19:03"
$ws.Range("K20").Value = "16.3 MJ/m2"
$ws.Range("E21").Value = "2026-02-21 17:19:06"
$ws.Range("J21").Value = "1030.8 hPa"
$ws.Range("K21").Value = "15.9 MJ/m2"
$ws.Range("O21").Value = "7.1 °C"
$ws.Range("E22").Value = "2026-02-21 17:19:08"
$ws.Range("K22").Value = "16.6 MJ/m2"
$ws.Range("E23").Value = "2026-02-21 17:19:11"
$ws.Range("K23").Value = "16.0 MJ/m2"
$ws.Range("O23").Value = "2.6 °C"
$ws.Range("E24").Value = "2026-02-21 17:19:13"
$ws.Range("J24").Value = "1031.7 hPa"
$ws.Range("K24").Value = "15.8 MJ/m2"
$ws.Range("O24").Value = "6.1 °C"
$ws.Range("E25").Value = "2026-02-21 17:19:16"
$ws.Range("K25").Value = "16.9 MJ/m2"
$ws.Range("O25").Value = "4.0 °C"
$ws.Range("E26").Value = "2026-02-21 17:19:18"
$ws.Range("E27").Value = "2026-02-21 17:19:21"
$ws.Range("K27").Value = "16.5 MJ/m2"
$ws.Range("E28").Value = "2026-02-21 17:19:23"
$ws.Range("K28").Value = "15.0 MJ/m2"
$ws.Range("O28").Value = "8.0 °C"
$ws.Range("E29").Value = "2026-02-21 17:19:26"
$ws.Range("E30").Value = "2026-02-21 17:19:28"
$ws.Range("K30").Value = "15.2 MJ/m2"
$ws.Range("E31").Value = "2026-02-21 17:19:30"
$ws.Range("K31").Value = "15.1 MJ/m2"
$ws.Range("O31").Value = "12.0 °C"
$ws.Range("E32").Value = "2026-02-21 17:19:33"
$ws.Range("K32").Value = "16.1 MJ/m2"
$ws.Range("O32").Value = "5.7 °C"
$ws.Range("E33").Value = "2026-02-21 17:19:35"
$ws.Range("J33").Value = "1030.5 hPa"
$ws.Range("K33").Value = "15.7 MJ/m2"
$ws.Range("O33").Value = "6.1 °C"
$ws.Range("E34").Value = "2026-02-21 17:19:38"
$ws.Range("K34").Value = "14.2 MJ/m2"
$ws.Range("E35").Value = "2026-02-21 17:19:40"
$ws.Range("J35").Value = "1030.7 hPa"
$ws.Range("K35").Value = "16.4 MJ/m2"
$ws.Range("O35").Value = "7.9 °C"
$ws.Range("E36").Value = "2026-02-21 17:19:42"
$ws.Range("K36").Value = "15.2 MJ/m2"
$ws.Range("E37").Value = "2026-02-21 17:19:45"
$ws.Range("O37").Value = "5.8 °C"
$ws.Range("E38").Value = "2026-02-21 17:19:47"
$ws.Range("K38").Value = "15.6 MJ/m2"
$ws.Range("O38").Value = "9.7 °C"
$ws.Range("E39").Value = "2026-02-21 17:19:50"
$ws.Range("E40").Value = "2026-02-21 17:19:52"
$ws.Range("J40").Value = "1030.4 hPa"
$ws.Range("O40").Value = "9.1 °C"
$ws.Range("E41").Value = "2026-02-21 17:19:54"
$ws.Range("K41").Value = "15.4 MJ/m2"
$ws.Range("O41").Value = "11.3 °C"
$ws.Range("E42").Value = "2026-02-21 17:19:57"
$ws.Range("O42").Value = "11.0 °C"
$ws.Range("E43").Value = "2026-02-21 17:19:59"
$ws.Range("K43").Value = "15.1 MJ/m2"
$ws.Range("O43").Value = "6.4 °C"
$ws.Range("E44").Value = "2026-02-21 17:20:02"
$ws.Range("K44").Value = "16.0 MJ/m2"
$ws.Range("E45").Value = "2026-02-21 17:20:04"
$ws.Range("J45").Value = "1032.3 hPa"
$ws.Range("O45").Value = "5.4 °C"
$ws.Range("E46").Value = "2026-02-21 17:20:07"
$ws.Range("J46").Value = "1031.7 hPa"
$ws.Range("K46").Value = "15.4 MJ/m2"
$ws.Range("O46").Value = "9.8 °C"

# Percentage cells must be written as literal text "NN%" rather than a numeric
# percentage, matching the inlineStr representation in the source file. Excel
# auto-converts a directly assigned "NN%" string into a numeric percent value, so
# instead we stage the text in a scratch cell below the table (forced to Text
# format), copy it, and paste-special (values only) into the target cell. This
# preserves the destination cell's existing style while keeping the value as text.
$scratch = $ws.Range("A47")
$scratch.NumberFormat = "@"
$scratch.Value = "42%"
$scratch.Copy()
$ws.Range("H2").PasteSpecial(-4163)
$scratch.Value = "69%"
$scratch.Copy()
$ws.Range("H4").PasteSpecial(-4163)
$scratch.Value = "70%"
$scratch.Copy()
$ws.Range("H6").PasteSpecial(-4163)
$scratch.Value = "50%"
$scratch.Copy()
$ws.Range("H9").PasteSpecial(-4163)
$scratch.Value = "77%"
$scratch.Copy()
$ws.Range("H10").PasteSpecial(-4163)
$scratch.Value = "48%"
$scratch.Copy()
$ws.Range("H11").PasteSpecial(-4163)
$scratch.Value = "56%"
$scratch.Copy()
$ws.Range("H12").PasteSpecial(-4163)
$scratch.Value = "62%"
$scratch.Copy()
$ws.Range("H13").PasteSpecial(-4163)
$scratch.Value = "34%"
$scratch.Copy()
$ws.Range("H16").PasteSpecial(-4163)
$scratch.Value = "74%"
$scratch.Copy()
$ws.Range("H18").PasteSpecial(-4163)
$scratch.Value = "65%"
$scratch.Copy()
$ws.Range("H19").PasteSpecial(-4163)
$scratch.Value = "57%"
$scratch.Copy()
$ws.Range("H21").PasteSpecial(-4163)
$scratch.Value = "33%"
$scratch.Copy()
$ws.Range("H23").PasteSpecial(-4163)
$scratch.Value = "83%"
$scratch.Copy()
$ws.Range("H24").PasteSpecial(-4163)
$scratch.Value = "72%"
$scratch.Copy()
$ws.Range("H28").PasteSpecial(-4163)
$scratch.Value = "55%"
$scratch.Copy()
$ws.Range("H33").PasteSpecial(-4163)
$scratch.Value = "57%"
$scratch.Copy()
$ws.Range("H35").PasteSpecial(-4163)
$scratch.Value = "78%"
$scratch.Copy()
$ws.Range("H43").PasteSpecial(-4163)
$scratch.Value = "68%"
$scratch.Copy()
$ws.Range("H46").PasteSpecial(-4163)

# Remove the scratch row entirely (with an upward shift) so the worksheet
# dimension and row contents return to their original bounds (A1:O46).
$ws.Rows("47:47").Delete()
